$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows (for line7 / line8) before the current row 8 ("extr1"),
# pushing the existing extr1..extr8 rows (old rows 8-15) down to rows 10-17.
$ws.Rows("8:9").Insert()

# Match the "A" column header style used by every other data row (bold,
# centered, thin box border) since the freshly-inserted rows don't inherit it.
$aStyle = $ws.Range("A8:A9")
$aStyle.Font.Bold = $true
$aStyle.HorizontalAlignment = -4108
$aStyle.VerticalAlignment = -4160
$aStyle.Borders.LineStyle = 1

# --- New rows 8 & 9: line7 / line8 contingency data ---
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# --- Rows 10-15 (formerly 8-13, extr1..extr6): refresh A/C/D/E values ---
$ws.Range("A10").Value = 8
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

$ws.Range("A11").Value = 9
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = $true

$ws.Range("A12").Value = 10
$ws.Range("C12").Value = 10
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = $false

$ws.Range("A13").Value = 11
$ws.Range("C13").Value = 7
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

$ws.Range("A14").Value = 12
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = $false

$ws.Range("A15").Value = 13
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $true

# --- New rows 16 & 17: extr7 / extr8 (previously rows 14-15) ---
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $false
